$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore the "min WCSS" values in column H (rows 3-24) that had been
# commented / cleared out, getting the sheet back to the state of Q1.
$values = @(
    22.530999999999999,
    9.7799999999999994,
    1.079,
    0.505,
    0.55400000000000005,
    0.49199999999999999,
    4.6669999999999998,
    0.47,
    1.417,
    0.40799999999999997,
    0.79700000000000004,
    1.218,
    0.45700000000000002,
    0.53400000000000003,
    0.40899999999999997,
    0.374,
    0.46400000000000002,
    0.625,
    0.34799999999999998,
    0.63100000000000001,
    2.1469999999999998,
    1.2430000000000001
)

$startRow = 3
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 8).Value = $values[$i]
}

# Page orientation explicitly set to portrait.
$ws.PageSetup.Orientation = 1

# Selection moved to the last filled cell in column H.
$ws.Range("H24").Select()
